$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Insert a new row below the "Model" row (row 8) to hold the new
# "L_curve" parameter, then relabel row 8 itself to "production_function".
$ws.Rows("9").Insert()

$ws.Range("A8").Value = "production_function"
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0

# Make optimization_parameters the active sheet/tab (was threshold_b)
# and select A9:B9 (the new L_curve row) there.
$ws.Activate()
$ws.Range("A9:B9").Select()
